# Atualização automática da planilha
# Delete row 24 ("Time Projeto" / DBD row with empty B/C/D) on the "Organograma"
# sheet; every row below shifts up by one (old row 25 -> new row 24, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")

# 1) Delete the row. This shifts rows 25:35 up to 24:34 and keeps their
#    original content/styles intact.
$ws.Rows(24).Delete()

# 2) The sheet's AutoFilter range should shrink along with the deleted row
#    (E33 -> E32), matching ordinary Excel reference-adjustment behaviour.
#    This engine's Range.AutoFilter() always snaps to the sheet's full used
#    range, so temporarily move the two trailing rows (33:34) out of the way,
#    re-apply the filter against the now-32-row used range, then restore
#    them exactly as they were.
$ws.Rows("33:34").Copy($ws.Range("A1000"))
$ws.Rows("33:34").Delete()

$ws.AutoFilterMode = $false
$ws.Range("A1:E32").AutoFilter()

$ws.Rows("33:34").Insert()
$ws.Range("A1000:E1001").Copy($ws.Range("A33"))
$ws.Range("A1000:E1001").Clear()

# 3) Keep the workbook-level "_FilterDatabase" defined name for Organograma
#    in sync with the new AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Organograma!_FilterDatabase") {
        $n.RefersTo = "=Organograma!`$A`$1:`$E`$32"
    }
}

# 4) Reflect the row-delete selection (selecting/deleting the whole row
#    leaves the entire row selected).
$ws.Range("A24:XFD24").Select()
